# "added pecota and more picks"
# Append new draft-pick rows to the "draftpicks" sheet (sheet1).
#
# NOTE: the row that ends up LAST (row 404, "Kevin Gausman") is entered
# first so that the shared-string table gets new unique strings appended
# in the same order the original author typed them in (Kevin Gausman was
# apparently keyed in before the rows that ultimately sort above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("draftpicks")

# Row 404 (entered first so "Kevin Gausman" becomes the first new shared string)
$ws.Cells.Item(404, 1).Value = "chicago"
$ws.Cells.Item(404, 2).Value = "Kevin Gausman"
$ws.Cells.Item(404, 3).Value = 3
$ws.Cells.Item(404, 4).Value = "P"
$ws.Cells.Item(404, 5).Value = 43867

# Row 397
$ws.Cells.Item(397, 1).Value = "bears"
$ws.Cells.Item(397, 2).Value = "Kevin Pillar"
$ws.Cells.Item(397, 3).Value = 1
$ws.Cells.Item(397, 4).Value = "OF"
$ws.Cells.Item(397, 5).Value = 43867

# Row 398
$ws.Cells.Item(398, 1).Value = "sturgeon"
$ws.Cells.Item(398, 2).Value = "Josh Lindblom"
$ws.Cells.Item(398, 3).Value = 4
$ws.Cells.Item(398, 4).Value = "P"
$ws.Cells.Item(398, 5).Value = 43867

# Row 399
$ws.Cells.Item(399, 1).Value = "balco"
$ws.Cells.Item(399, 2).Value = "Austin Hedges"
$ws.Cells.Item(399, 3).Value = 1
$ws.Cells.Item(399, 4).Value = "C"
$ws.Cells.Item(399, 5).Value = 43867

# Row 400
$ws.Cells.Item(400, 1).Value = "chicago"
$ws.Cells.Item(400, 2).Value = "Garrett Richards"
$ws.Cells.Item(400, 3).Value = 6
$ws.Cells.Item(400, 4).Value = "P"
$ws.Cells.Item(400, 5).Value = 43867

# Row 401
$ws.Cells.Item(401, 1).Value = "dembums"
$ws.Cells.Item(401, 2).Value = "Mike Tauchman"
$ws.Cells.Item(401, 3).Value = 1
$ws.Cells.Item(401, 4).Value = "OF"
$ws.Cells.Item(401, 5).Value = 43867

# Row 402
$ws.Cells.Item(402, 1).Value = "ds9"
$ws.Cells.Item(402, 2).Value = "Shogo Akiyama"
$ws.Cells.Item(402, 3).Value = 5
$ws.Cells.Item(402, 4).Value = "OF"
$ws.Cells.Item(402, 5).Value = 43867

# Row 403
$ws.Cells.Item(403, 1).Value = "drjames"
$ws.Cells.Item(403, 2).Value = "Stephen Piscotty"
$ws.Cells.Item(403, 3).Value = 3
$ws.Cells.Item(403, 4).Value = "OF"
$ws.Cells.Item(403, 5).Value = 43867

# Match the final cursor position left in the saved workbook
$ws.Range("D404").Select()
